# Apply the "Update CDA Logical model for ST.r2b" edit:
#  - rename the include sheet
#  - bump Version / Date metadata values
#  - insert a new "Jurisdiction" metadata row (with a blank value) right
#    after "Contact", pushing Description/Purpose/Copyright/Immutable down

$wb = $excel.ActiveWorkbook

$metaSheetName = "Metadata"
$includeSheetName = "Include from EntityNameUse"

$ws1 = $wb.Worksheets.Item($metaSheetName)
$ws2 = $wb.Worksheets.Item($includeSheetName)

# 1) Rename the "Include from EntityNameUse" sheet to "Include #0"
$ws2.Name = "Include #0"

# 2) Update the Version value (row 3) and Date value (row 8)
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# 3) Insert a new row for "Jurisdiction" after "Contact" (row 10), before
#    "Description" (old row 11) - matches the style of the surrounding rows
$ws1.Rows.Item(11).Insert()
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""
